# Split the table: insert two new columns ("prop1" and "prop2" series)
# right after column A ("members"), pushing the existing header/data
# columns two places to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at B:C - everything currently in B onward
# shifts right by two columns.
$ws.Range("B1:C1").EntireColumn.Insert()

# Fill the new "prop1" column (B) - header + 11 members' values.
$ws.Range("B1").Value = "prop1"
$ws.Range("B2").Value = "prop1-01"
$ws.Range("B3").Value = "prop1-02"
$ws.Range("B4").Value = "prop1-03"
$ws.Range("B5").Value = "prop1-04"
$ws.Range("B6").Value = "prop1-05"
$ws.Range("B7").Value = "prop1-06"
$ws.Range("B8").Value = "prop1-07"
$ws.Range("B9").Value = "prop1-08"
$ws.Range("B10").Value = "prop1-09"
$ws.Range("B11").Value = "prop1-10"
$ws.Range("B12").Value = "prop1-11"

# Fill the new "prop2" column (C) - header + 11 members' values.
$ws.Range("C1").Value = "prop2"
$ws.Range("C2").Value = "prop2-01"
$ws.Range("C3").Value = "prop2-02"
$ws.Range("C4").Value = "prop2-03"
$ws.Range("C5").Value = "prop2-04"
$ws.Range("C6").Value = "prop2-05"
$ws.Range("C7").Value = "prop2-06"
$ws.Range("C8").Value = "prop2-07"
$ws.Range("C9").Value = "prop2-08"
$ws.Range("C10").Value = "prop2-09"
$ws.Range("C11").Value = "prop2-10"
$ws.Range("C12").Value = "prop2-11"

# Keep the active selection consistent with the edited workbook (F12,
# matching the post-edit cursor position recorded in the sheet view).
$ws.Range("F12").Select()
